$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire row 429 ("「ママと一緒にお菓子作り」...") which shifts
# every subsequent row (430-541) up by one, turning them into rows 429-540.
$ws.Rows("429:429").Delete()
